# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.089.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.67%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.732.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.37%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.14%  "

# Row 6
$ws.Range("E6").Value = "  -0.05%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4873"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.54%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3519"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.85%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.52"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.99%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07291"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.81%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.053"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.80%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.09%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.57%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.899"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.23%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.732.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.44%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.914"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.53%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.43%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001041"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.21%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06413"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.05%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.06%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.97%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.713"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.70%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.134.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.64%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.80%  "

# Row 25
$ws.Range("E25").Value = "  -3.59%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.74%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.18%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.931.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.48%  "

# Row 29
$ws.Range("E29").Value = "  -3.20%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.44%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.054"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.48%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09350"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.07%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.643"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.38%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.401"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.36%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05975"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.00%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02195"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.02%  "

# Row 37
$ws.Range("B37").Value = "WEMIXTOKEN"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.440"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.35%  "

# Row 38
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "11.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.10%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.785"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.15%  "

# Row 40
$ws.Range("E40").Value = "  -2.79%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6020"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.52%  "

# Row 42
$ws.Range("E42").Value = "  +0.14%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.098"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.86%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.524"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.00%  "

# Row 45
$ws.Range("E45").Value = "  -2.28%  "

# Row 46
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.587"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.75%  "

# Row 47
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5680"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.93%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "119.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.80%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.858"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.55%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.105"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.37%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06654"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.01%  "
